# "Fix kickoff after PAT"
# The "Ongoing Games" sheet tracks a live play-by-play snapshot of each
# game. This updates the snapshot row for the Omaha @ James Madison game
# (row 3) to reflect the next play after a PAT: James Madison now has the
# ball (kickoff return), it's a NORMAL play (not a KICKOFF) and the
# waiting-on user/score/clock/down fields move along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ongoing Games")

# Home score ticks up (extra point / score change reflected on the sheet)
$ws.Range("G3").Value = 7

# Game clock advances
$ws.Range("AC3").Value = "0:37"

# Ball spot moves
$ws.Range("AD3").Value = "James Madison 20"

# Down resets
$ws.Range("AE3").Value = 2

# Possession flips to James Madison after the kickoff
$ws.Range("AG3").Value = "James Madison"

# Now waiting on the James Madison user
$ws.Range("AH3").Value = "door_nav#2953"

# No longer waiting on a kickoff - next play is a normal snap
$ws.Range("AI3").Value = "NORMAL"

# Offensive / defensive random numbers rolled for the next play
$ws.Range("AJ3").Value = 664
$ws.Range("AK3").Value = 1212

# Number no longer submitted for the new play
$ws.Range("AP3").Value = "NO"

# Reflect the updated selection/cursor position left in the sheet by the
# editor (one column further right than before, matching the "Waiting On"
# column that was just edited).
$ws.Range("AH3").Select() | Out-Null
